$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.861.24'
$ws.Range('D3').Value = '2.533.48'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.13'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.74'
$ws.Range('E6').Value = '  +1.17%  '
$ws.Range('E7').Value = '  -1.12%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.537'
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.85'
$ws.Range('E10').Value = '  -1.53%  '
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.58'
$ws.Range('E12').Value = '  -1.18%  '
$ws.Range('E13').Value = '  -3.60%  '
$ws.Range('D14').Value = '2.921.68'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.500.15'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.15'
$ws.Range('E16').Value = '  -4.01%  '
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '42.909.69'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.88'
$ws.Range('E19').Value = '  +2.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.76'
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('D21').Value = '0.0₃0968'
$ws.Range('E21').Value = '  -0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.66'
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '253.16'
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.07'
$ws.Range('E25').Value = '  +1.31%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.31'
$ws.Range('E26').Value = '  -4.32%  '
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('E28').Value = '  +2.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.88'
$ws.Range('E29').Value = '  +2.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.53'
$ws.Range('E30').Value = '  +4.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.90'
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.38'
$ws.Range('E32').Value = '  +1.10%  '
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.30'
$ws.Range('E34').Value = '  -3.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.71'
$ws.Range('E35').Value = '  +3.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.34'
$ws.Range('E36').Value = '  -2.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0793'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.46'
$ws.Range('E39').Value = '  +8.43%  '
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('E41').Value = '  -11.23%  '
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.28'
$ws.Range('E45').Value = '  -3.58%  '
$ws.Range('D46').Value = '2.012.25'
$ws.Range('E46').Value = '  -1.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.19'
$ws.Range('E47').Value = '  +3.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '84.19'
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '106.62'
$ws.Range('E49').Value = '  +4.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.00'
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('D51').Value = '2.776.87'
$ws.Range('E51').Value = '  +0.13%  '
